# Switch the presentation's design theme from "Integral" (Red Violet colour
# scheme) to "Office Theme" (Office colour scheme).
#
# PowerPoint models a deck's colour scheme as the 12-slot ThemeColorScheme
# hanging off the (single) slide master's Theme object:
#   dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
# RGB values are plain OLE_COLOR longs (0xBBGGRR).

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Office Theme colour values (srgb hex -> OLE BGR long)
$officeColors = @(
    0x000000,  # dk1      000000
    0xFFFFFF,  # lt1      FFFFFF
    0x6A5444,  # dk2      44546A -> BGR
    0xE6E6E7,  # lt2      E7E6E6 -> BGR
    0xD59B5B,  # accent1  5B9BD5 -> BGR
    0x317DED,  # accent2  ED7D31 -> BGR
    0xA5A5A5,  # accent3  A5A5A5 -> BGR
    0x00C0FF,  # accent4  FFC000 -> BGR
    0xC47244,  # accent5  4472C4 -> BGR
    0x47AD70,  # accent6  70AD47 -> BGR
    0xC16305,  # hlink    0563C1 -> BGR
    0x724F95   # folHlink 954F72 -> BGR
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
